# 06/11/2025 Fixed Showing up of Draft Button for IT PIC (Accepted Tix)
#
# Adds six new report columns (P:U) to the "ResolvedTicket" sheet header
# row (row 4): ASSIGNED IT PIC, ASSIGNED DATE TIME, RESOLVED DATE TIME,
# SLA HOURS, ACTUAL HOURS and HIT OR MISS - plus matching column widths,
# the header cell fill/shading style, and refreshed view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the header shading/fill style from the last existing header
# cell (O4) onto the new header cells before filling them in, so P4:U4
# pick up the same "s=1" (shaded fill) cell style used by A4:O4.
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4:U4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Fill in the new header labels. The order below (not strict
# left-to-right column order) matches how the shared-string table was
# populated in the authored workbook.
$ws.Range("Q4").Value = "ASSIGNED DATE TIME"
$ws.Range("P4").Value = "ASSIGNED IT PIC"
$ws.Range("U4").Value = "HIT OR MISS"
$ws.Range("R4").Value = "RESOLVED DATE TIME"
$ws.Range("S4").Value = "SLA HOURS"
$ws.Range("T4").Value = "ACTUAL HOURS"

# --- Match the new columns' widths to the target layout. ColumnWidth is
# character-width based and gets pixel-snapped on save (same as real
# Excel), so these are the closest achievable settings.
$ws.Columns.Item(16).ColumnWidth = 32.83333333333333  # P: ASSIGNED IT PIC
$ws.Columns.Item(17).ColumnWidth = 36.83333333333333  # Q: ASSIGNED DATE TIME
$ws.Columns.Item(18).ColumnWidth = 45.0                # R: RESOLVED DATE TIME
$ws.Columns.Item(19).ColumnWidth = 20.666666666666664  # S: SLA HOURS
$ws.Columns.Item(20).ColumnWidth = 23.166666666666664  # T: ACTUAL HOURS
$ws.Columns.Item(21).ColumnWidth = 18.666666666666664  # U: HIT OR MISS

# --- Refresh the view state: zoom to 85% and move the selection to
# reflect the area the new columns were authored in.
$excel.ActiveWindow.Zoom = 85
$ws.Range("S14").Select() | Out-Null
